$d = $word.ActiveDocument

# --- Create the three new character styles (appended to styles.xml) ---

$gaNStyle = $d.Styles.Add("GaNStyle", 2)
$gaNStyle.Font.Name = "Calibri"
$gaNStyle.Font.Size = 14

$gaNParagraph = $d.Styles.Add("GaNParagraph", 2)
$gaNParagraph.Font.Name = "Calibri"
$gaNParagraph.Font.Size = 10

$gaNLinks = $d.Styles.Add("GaNLinks", 2)
$gaNLinks.Font.Name = "Calibri"
$gaNLinks.Font.Bold = $true
$gaNLinks.Font.Color = 8388608
$gaNLinks.Font.Size = 9.5
$gaNLinks.Font.Underline = 1

# --- Apply GaNStyle to every "2022 observing dates" heading run (4 occurrences) ---

$headingText = "2022 Ημερομηνίες παρατήρησης για τον  Αστερισμός του Κύκνου: 10-19 Αυγούστου, 9-18 Σεπτεμβρίου, 8-17 Οκτωβρίου"
$rng = $d.Content
while ($rng.Find.Execute($headingText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Style = "GaNStyle"
    $rng.Collapse(0)
}

# --- Apply GaNParagraph to the campaign description paragraph ---

$paraText = "Συμμετέχετε σε μία παγκόσμια καμπάνια για να παρατηρήσετε και να καταγράψετε τη φωτεινότητα των πιο αμυδρά ορατών άστρων σαν μέσο για την μέτρηση της Φωτορρύπανσης σε μία δεδομένη περιοχή. Με τον εντοπισμό και την παρατήρηση του  Αστερισμός του Κύκνου στον νυχτερινό ουρανό καθώς και με την σύγκριση των ανωτέρω με τα διαγράμματα για τα μεγέθη των άστρων,  άνθρωποι από όλον τον κόσμο θα μάθουν πώς τα φώτα στην κοινότητά τους συμβάλλουν στην Φωτορρύπανση. Με την κατάθεση των πορισμάτων τους στην ιστοσελίδα θα δημιουργηθεί ένα αρχείο σχετικά με το τι μπορεί να δει κανείς στον νυχτερινό ουρανό."
$rng2 = $d.Content
if ($rng2.Find.Execute($paraText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng2.Style = "GaNParagraph"
}

# --- Apply GaNLinks to the chart-credit paragraph ---

$linksText = "Τα διαγράμματα αυτού του αρχείου επιμελήθηκε ο Jan Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."
$rng3 = $d.Content
if ($rng3.Find.Execute($linksText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng3.Style = "GaNLinks"
}

Write-Output "done"
